$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("8.00", "0.0190", etc.)
# must be forced to Text format first, otherwise Excel auto-converts them
# (stripping the trailing zeros / thousands-style dots) the same way typing
# such text into a General-formatted cell would.
$textCells = @("D5", "D6", "D9", "D13", "D15", "D16", "D18", "D20", "D21", "D24", "D26", "D27", "D31", "D32", "D36", "D38", "D39", "D40", "D41", "D44", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.530.76"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.813.24"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "228.73"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "0.577"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").Value = "0.302"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D12").Value = "2.073.89"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "11.26"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "1.814.73"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "0.648"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.47"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "34.511.94"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "69.22"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "245.65"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "11.47"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").Value = "172.37"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  +9.01%  "
$ws.Range("D27").Value = "16.84"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "0.0534"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").Value = "3.87"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "1.396.70"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("D36").Value = "0.682"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "1.07"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "83.94"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").Value = "0.965"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "13.37"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").Value = "0.0514"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "6.00"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "1.973.48"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "105.30"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +0.20%  "
